# Apply updates from EMH project for the Ref_LDV ZEV BC reference policy workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reword the comment in X3: move the "2015 & 2020 values are actual sales"
# clause earlier in the sentence (parenthetical), and tidy punctuation.
$ws.Range("X3").Value = "Should move 2015 & 2020 values (actual sales) once include subsidy policy! ; https://www150.statcan.gc.ca/t1/tbl1/en/cv.action?pid=2010002101"

# ZEV market share_class_min values for 2035-2050 (T3:W3) bumped from 0.99 to 1 (100%).
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 1
